$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 91

# Columns A (date-looking "2025-12-16") and C (numeric-looking "251216")
# would otherwise be auto-coerced by Excel into a real date / number.
# Force them to text first so the literal strings are preserved, exactly
# like the other rows in this sheet (all values are stored as text).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-12-16"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251216"
$ws.Cells.Item($newRow, 4).Value = "8-9-2"
$ws.Cells.Item($newRow, 5).Value = "2025-12-16T21:46:07.770+04:00"
